# Applies the coin price/volume update described in the commit:
# "Updated symbol list on Wed Feb 15 11:27:51 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (leading "'" forces Excel to
# store it as a string instead of re-parsing "42.21"/"2.35%"/etc. as a
# number/percentage) while keeping the General number format intact,
# matching the workbook's existing inlineStr text cells.
function Set-TextCell {
    param($Sheet, $Addr, $Val)
    $Sheet.Range($Addr).Value = "'" + $Val
}

Set-TextCell $ws 'E2' '2.35%'
Set-TextCell $ws 'D3' '42.21'
Set-TextCell $ws 'E3' '4.53%'
Set-TextCell $ws 'D4' '5.014'
Set-TextCell $ws 'E4' '0.17%'
Set-TextCell $ws 'D5' '0.07570'
Set-TextCell $ws 'E5' '3.32%'
Set-TextCell $ws 'B6' 'FTXToken'
Set-TextCell $ws 'C6' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell $ws 'D6' '1.600'
Set-TextCell $ws 'E6' '2.38%'
Set-TextCell $ws 'B7' 'MXToken'
Set-TextCell $ws 'C7' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws 'D7' '0.9412'
Set-TextCell $ws 'E7' '1.94%'
Set-TextCell $ws 'B8' 'BTSEToken'
Set-TextCell $ws 'C8' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell $ws 'D8' '2.384'
Set-TextCell $ws 'E8' '0.14%'
Set-TextCell $ws 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextCell $ws 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell $ws 'D9' '0.1196'
Set-TextCell $ws 'E9' '1.32%'
Set-TextCell $ws 'B10' 'WazirX'
Set-TextCell $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell $ws 'D10' '0.1837'
Set-TextCell $ws 'E10' '1.12%'
Set-TextCell $ws 'B11' 'MandalaExchangeToken'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell $ws 'D11' '0.09129'
Set-TextCell $ws 'E11' '3.28%'
Set-TextCell $ws 'B12' 'BitrueCoin'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell $ws 'D12' '0.04198'
Set-TextCell $ws 'E12' '-4.39%'
Set-TextCell $ws 'B13' 'BitMartToken'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell $ws 'D13' '0.1049'
Set-TextCell $ws 'E13' '-0.38%'
Set-TextCell $ws 'B14' 'BitForexToken'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell $ws 'D14' '0.001287'
Set-TextCell $ws 'E14' '1.91%'
Set-TextCell $ws 'B15' 'TigerCash'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws 'D15' '0.005782'
Set-TextCell $ws 'E15' '-1.13%'
Set-TextCell $ws 'B16' 'LEO'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws 'D16' '3.344'
Set-TextCell $ws 'E16' '0.11%'
Set-TextCell $ws 'B17' 'GateToken'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell $ws 'D17' '4.382'
Set-TextCell $ws 'E17' '2.17%'
Set-TextCell $ws 'E18' '0.21%'
Set-TextCell $ws 'D19' '8.382'
Set-TextCell $ws 'E19' '6.52%'
Set-TextCell $ws 'E20' '1.36%'
Set-TextCell $ws 'D21' '0.3298'
Set-TextCell $ws 'E21' '16.01%'
Set-TextCell $ws 'D22' '0.04104'
Set-TextCell $ws 'E22' '4.78%'
Set-TextCell $ws 'D23' '0.001265'
Set-TextCell $ws 'E23' '0.34%'
Set-TextCell $ws 'D24' '0.003895'
Set-TextCell $ws 'E24' '5.19%'
Set-TextCell $ws 'E25' '1.41%'
Set-TextCell $ws 'D38' '0.02411'
Set-TextCell $ws 'E38' '2.94%'
Set-TextCell $ws 'D39' '0.05248'
Set-TextCell $ws 'E39' '3.25%'
Set-TextCell $ws 'D40' '0.006676'
Set-TextCell $ws 'E40' '10.94%'
Set-TextCell $ws 'D41' '0.007698'
Set-TextCell $ws 'D42' '0.1330'
Set-TextCell $ws 'E42' '3.06%'
Set-TextCell $ws 'D43' '0.007389'
Set-TextCell $ws 'E43' '0.01%'
Set-TextCell $ws 'D44' '0.007801'
Set-TextCell $ws 'E44' '-3.24%'
Set-TextCell $ws 'D45' '0.3005'
Set-TextCell $ws 'E45' '3.37%'
Set-TextCell $ws 'D46' '0.00006243'
Set-TextCell $ws 'E46' '0.60%'
Set-TextCell $ws 'D47' '0.00000000750'
Set-TextCell $ws 'E47' '-0.17%'
Set-TextCell $ws 'D48' '0.04524'
Set-TextCell $ws 'E48' '-4.28%'
Set-TextCell $ws 'E49' '-0.05%'
Set-TextCell $ws 'D50' '0.00002099'
Set-TextCell $ws 'E50' '-0.17%'
Set-TextCell $ws 'E51' '-0.17%'

Write-Output "Applied 93 cell updates"